$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting of the last existing data row (328) down through the new rows (329:343)
$ws.Range("A328:D328").Copy() | Out-Null
$ws.Range("A329:D343").PasteSpecial(-4122) | Out-Null

# Fill in the new data values
$ws.Range("A329").Value = 44403
$ws.Range("B329").Value = 6
$ws.Range("C329").Value = 90
$ws.Range("D329").Value = 47.61577246009534

$ws.Range("A330").Value = 44404
$ws.Range("B330").Value = 17
$ws.Range("C330").Value = 103
$ws.Range("D330").Value = 54.49360625988689

$ws.Range("A331").Value = 44405
$ws.Range("B331").Value = 12
$ws.Range("C331").Value = 109
$ws.Range("D331").Value = 57.66799109055991

$ws.Range("A332").Value = 44406
$ws.Range("B332").Value = 22
$ws.Range("C332").Value = 111
$ws.Range("D332").Value = 58.72611936745092

$ws.Range("A333").Value = 44407
$ws.Range("B333").Value = 25
$ws.Range("C333").Value = 120
$ws.Range("D333").Value = 63.48769661346046

$ws.Range("A334").Value = 44408
$ws.Range("B334").Value = 7
$ws.Range("C334").Value = 111
$ws.Range("D334").Value = 58.72611936745092

$ws.Range("A335").Value = 44409
$ws.Range("B335").Value = 32
$ws.Range("C335").Value = 121
$ws.Range("D335").Value = 64.01676075190595

$ws.Range("A336").Value = 44410
$ws.Range("B336").Value = 29
$ws.Range("C336").Value = 144
$ws.Range("D336").Value = 76.18523593615254

$ws.Range("A337").Value = 44411
$ws.Range("B337").Value = 4
$ws.Range("C337").Value = 131
$ws.Range("D337").Value = 69.30740213636099

$ws.Range("A338").Value = 44412
$ws.Range("B338").Value = 20
$ws.Range("C338").Value = 139
$ws.Range("D338").Value = 73.53991524392502

$ws.Range("A339").Value = 44413
$ws.Range("B339").Value = 29
$ws.Range("C339").Value = 146
$ws.Range("D339").Value = 77.24336421304355

$ws.Range("A340").Value = 44414
$ws.Range("B340").Value = 17
$ws.Range("C340").Value = 138
$ws.Range("D340").Value = 73.01085110547952

$ws.Range("A341").Value = 44415
$ws.Range("B341").Value = 13
$ws.Range("C341").Value = 144
$ws.Range("D341").Value = 76.18523593615254

$ws.Range("A342").Value = 44416
$ws.Range("B342").Value = 24
$ws.Range("C342").Value = 136
$ws.Range("D342").Value = 71.9527228285885

$ws.Range("A343").Value = 44417
$ws.Range("B343").Value = 16
$ws.Range("C343").Value = 123
$ws.Range("D343").Value = 65.07488902879696

